$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Cells.Item(2, 4).Value = '27.024.57'
$ws.Cells.Item(2, 5).Value = '  -0.19%  '

# Row 3: Ethereum -> Ethereum
$ws.Cells.Item(3, 4).Value = '1.823.01'
$ws.Cells.Item(3, 5).Value = '  +1.16%  '

# Row 4: TetherUSD -> TetherUSD
$ws.Cells.Item(4, 4).Value = "'1.008"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5: BNB -> BNB
$ws.Cells.Item(5, 4).Value = "'312.20"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.33%  '

# Row 6: USDC -> USDC
$ws.Cells.Item(6, 4).Value = "'1.007"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.02%  '

# Row 7: XRP -> XRP
$ws.Cells.Item(7, 4).Value = "'0.4313"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +1.84%  '

# Row 8: Cardano -> Cardano
$ws.Cells.Item(8, 4).Value = "'0.3705"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +2.69%  '

# Row 9: Dogecoin -> Dogecoin
$ws.Cells.Item(9, 4).Value = "'0.07260"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +1.09%  '

# Row 10: WrappedEther -> Polygon
$ws.Cells.Item(10, 2).Value = 'Polygon'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10, 4).Value = "'0.8659"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +2.73%  '

# Row 11: Polygon -> WrappedEther
$ws.Cells.Item(11, 2).Value = 'WrappedEther'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(11, 4).Value = '2.054.73'
$ws.Cells.Item(11, 5).Value = '  +12.26%  '

# Row 12: Solana -> Solana
$ws.Cells.Item(12, 4).Value = "'21.02"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +3.88%  '

# Row 13: Chainlink -> Chainlink
$ws.Cells.Item(13, 4).Value = "'6.642"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +4.05%  '

# Row 14: Polkadot -> Polkadot
$ws.Cells.Item(14, 4).Value = "'5.410"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.45%  '

# Row 15: TRON -> TRON
$ws.Cells.Item(15, 4).Value = "'0.06933"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.70%  '

# Row 16: Litecoin -> Litecoin
$ws.Cells.Item(16, 4).Value = "'81.05"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.09%  '

# Row 17: BinanceUSD -> BinanceUSD
$ws.Cells.Item(17, 4).Value = "'1.012"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.10%  '

# Row 18: ShibaInu -> ShibaInu
$ws.Cells.Item(18, 4).Value = "'0.000008883"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +2.20%  '

# Row 19: Dai -> Dai
$ws.Cells.Item(19, 4).Value = "'1.007"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.06%  '

# Row 20: Avalanche -> Avalanche
$ws.Cells.Item(20, 4).Value = "'15.24"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.84%  '

# Row 21: WrappedBTC -> WrappedBTC
$ws.Cells.Item(21, 4).Value = '27.084.05'
$ws.Cells.Item(21, 5).Value = '  -0.40%  '

# Row 22: Uniswap -> Uniswap
$ws.Cells.Item(22, 4).Value = "'5.203"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.83%  '

# Row 23: Cosmos -> Cosmos
$ws.Cells.Item(23, 4).Value = "'11.12"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.20%  '

# Row 24: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Cells.Item(24, 4).Value = '2.298.92'
$ws.Cells.Item(24, 5).Value = '  +12.11%  '

# Row 25: Monero -> Monero
$ws.Cells.Item(25, 4).Value = "'154.34"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.01%  '

# Row 26: Toncoin -> Toncoin
$ws.Cells.Item(26, 4).Value = "'1.899"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.18%  '

# Row 27: EthereumClassic -> EthereumClassic
$ws.Cells.Item(27, 4).Value = "'18.34"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.80%  '

# Row 28: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Cells.Item(28, 4).Value = "'5.226"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +4.61%  '

# Row 29: LidoDAOToken -> LidoDAOToken
$ws.Cells.Item(29, 4).Value = "'1.903"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +17.93%  '

# Row 30: BitcoinCash -> BitcoinCash
$ws.Cells.Item(30, 4).Value = "'115.20"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.30%  '

# Row 31: Stellar -> Stellar
$ws.Cells.Item(31, 4).Value = "'0.08955"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.22%  '

# Row 32: ImmutableX -> ImmutableX
$ws.Cells.Item(32, 4).Value = "'0.7470"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +2.95%  '

# Row 33: ARBITRUM -> ARBITRUM
$ws.Cells.Item(33, 5).Value = '  +6.87%  '

# Row 34: Filecoin -> Filecoin
$ws.Cells.Item(34, 4).Value = "'4.427"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.16%  '

# Row 35: HuobiToken -> HuobiToken
$ws.Cells.Item(35, 4).Value = "'2.810"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.33%  '

# Row 36: Frax -> Frax
$ws.Cells.Item(36, 4).Value = "'1.013"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.52%  '

# Row 37: TrustWalletToken -> TrustWalletToken
$ws.Cells.Item(37, 4).Value = "'1.128"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +3.63%  '

# Row 38: Hedera -> Hedera
$ws.Cells.Item(38, 4).Value = "'0.05222"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +2.39%  '

# Row 39: VeChain -> VeChain
$ws.Cells.Item(39, 4).Value = "'0.01935"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +1.69%  '

# Row 40: TheSandbox -> TheSandbox
$ws.Cells.Item(40, 4).Value = "'0.5096"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +2.79%  '

# Row 41: Algorand -> MXToken
$ws.Cells.Item(41, 2).Value = 'MXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(41, 4).Value = "'2.757"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +5.76%  '

# Row 42: MXToken -> Algorand
$ws.Cells.Item(42, 2).Value = 'Algorand'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(42, 4).Value = "'0.1649"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +1.90%  '

# Row 43: FraxShare -> FraxShare
$ws.Cells.Item(43, 4).Value = "'6.502"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +8.64%  '

# Row 44: Aptos -> Aptos
$ws.Cells.Item(44, 4).Value = "'8.247"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +2.46%  '

# Row 45: Quant -> EnergySwap
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).Value = "'10.45"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.75%  '

# Row 46: EnergySwap -> Quant
$ws.Cells.Item(46, 2).Value = 'Quant'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(46, 4).Value = "'106.79"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.98%  '

# Row 47: PaxDollar -> PaxDollar
$ws.Cells.Item(47, 4).Value = "'1.008"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.01%  '

# Row 48: Cronos -> NEARProtocol
$ws.Cells.Item(48, 2).Value = 'NEARProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(48, 4).Value = "'1.655"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +4.42%  '

# Row 49: NEARProtocol -> Cronos
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).Value = "'0.06318"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.14%  '

# Row 50: Decentraland -> Decentraland
$ws.Cells.Item(50, 4).Value = "'0.4570"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +1.41%  '

# Row 51: RenderToken -> RenderToken
$ws.Cells.Item(51, 4).Value = "'1.802"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +6.02%  '
